# Updates the cryptos price/volume snapshot (GitHub Actions data refresh).
# Price cells (column D) that look like plain numbers get a leading "'"
# so Excel stores them as text (matching the original inlineStr cells)
# instead of silently parsing them into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.073.09"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "3.513.71"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'610.69"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "'148.27"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("D7").Value = "3.513.16"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.478"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").Value = "'0.142"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "'8.08"
$ws.Range("E11").Value = "  +7.33%  "
$ws.Range("D12").Value = "'0.422"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "'0.0000216"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "4.111.60"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "'31.58"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "3.520.30"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "67.181.30"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'10.84"
$ws.Range("E19").Value = "  +9.12%  "
$ws.Range("D20").Value = "'6.37"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "'15.36"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'436.79"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").Value = "'0.608"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").Value = "'79.62"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").Value = "3.653.52"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'0.0000120"
$ws.Range("E27").Value = "  -3.26%  "
$ws.Range("D28").Value = "'9.78"
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("D29").Value = "'8.21"
$ws.Range("E29").Value = "  -5.11%  "
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").Value = "'1.58"
$ws.Range("E31").Value = "  -4.76%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("D34").Value = "'25.57"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "'5.97"
$ws.Range("E35").Value = "  -3.16%  "
$ws.Range("D36").Value = "'1.80"
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").Value = "'8.02"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").Value = "'176.26"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "'0.0895"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "'5.41"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'2.04"
$ws.Range("E43").Value = "  -11.62%  "
$ws.Range("D44").Value = "'0.896"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "'46.28"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").Value = "'27.85"
$ws.Range("E46").Value = "  -6.61%  "
$ws.Range("D47").Value = "'1.25"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'7.47"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'2.45"
$ws.Range("E49").Value = "  -2.48%  "
$ws.Range("D50").Value = "'0.996"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").Value = "'0.246"
$ws.Range("E51").Value = "  -2.30%  "
